$d = $word.ActiveDocument

# The "<id>...</id>" tag text had been split across three runs (the
# Courier-New-styled "<id>" / "</id>" delimiters plus a separately
# formatted run for the inner value). Collapse each occurrence back
# into a single run carrying the tag's Courier New styling, for both
# places this pattern occurs in the document (p131r_5 and p131v_1).

$r1 = $d.Content
$r1.Find.Execute("<id>p131r_5</id>", $true, $false, $false, $false, $false, `
                  $false, 1, $false, "<id>p131r_5</id>", 2)

$r2 = $d.Content
$r2.Find.Execute("<id>p131v_1</id>", $true, $false, $false, $false, $false, `
                  $false, 1, $false, "<id>p131v_1</id>", 2)

Write-Output "done"
